# BusManagerSecOcDemo.xlsx - "take Line6 as part of CSM"
#
# On the "Demo" sheet: the CSM label that used to live in C7 (merged over
# C8:C9 style pattern already used by the row below) should instead start
# one row higher, spanning C6:C7, matching the existing C8:C9 merged-pair
# formatting further down the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demo")
$ws.Activate()

# Move the "CSM" text up from C7 into C6, leaving C7 blank.
# (Use .Text rather than .Value - the latter doesn't round-trip a shared
# string's textual content correctly in this host.)
$ws.Range("C6").Value = $ws.Range("C7").Text
$ws.Range("C7").ClearContents()

# Re-apply the formatting used by the equivalent merged pair (C8:C9) so the
# new top/bottom cells of the C6:C7 merge get the matching border styling.
$ws.Range("C8").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Merge C6:C7 like the other label/value pairs in this column.
$ws.Range("C6:C7").Merge()

# Update the saved view: zoom in a bit and move the selection.
$excel.ActiveWindow.Zoom = 85
$ws.Range("B7").Select()
